$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Clcf1"
$ws.Cells.Item(2, 3).Value = "Lifr"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.6103310000000001
$ws.Cells.Item(2, 8).Value = 1.830993
$ws.Cells.Item(2, 9).Value = 0.06059696746595619
$ws.Cells.Item(2, 10).Value = 0.06059696746595618
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 31.44605333333334
$ws.Cells.Item(2, 14).Value = 94.33816
$ws.Cells.Item(2, 15).Value = 0.273208187120734
$ws.Cells.Item(2, 16).Value = 0.273208187120734
$ws.Cells.Item(2, 17).Value = 19.19250117698667
$ws.Cells.Item(2, 18).Value = 172.73251059288
$ws.Cells.Item(2, 19).Value = 0.01655558762638799
$ws.Cells.Item(2, 20).Value = 0.01655558762638799

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Clcf1"
$ws.Cells.Item(3, 3).Value = "Lifr"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.6103310000000001
$ws.Cells.Item(3, 8).Value = 1.830993
$ws.Cells.Item(3, 9).Value = 0.06059696746595619
$ws.Cells.Item(3, 10).Value = 0.06059696746595618
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 57.80064033333333
$ws.Cells.Item(3, 14).Value = 173.401921
$ws.Cells.Item(3, 15).Value = 0.5021809252974907
$ws.Cells.Item(3, 16).Value = 0.5021809252974908
$ws.Cells.Item(3, 17).Value = 35.27752261528367
$ws.Cells.Item(3, 18).Value = 317.497703537553
$ws.Cells.Item(3, 19).Value = 0.03043064119227582
$ws.Cells.Item(3, 20).Value = 0.03043064119227583

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Clcf1"
$ws.Cells.Item(4, 3).Value = "Lifr"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.6103310000000001
$ws.Cells.Item(4, 8).Value = 1.830993
$ws.Cells.Item(4, 9).Value = 0.06059696746595619
$ws.Cells.Item(4, 10).Value = 0.06059696746595618
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 12.04190666666667
$ws.Cells.Item(4, 14).Value = 36.12572
$ws.Cells.Item(4, 15).Value = 0.1046219522368387
$ws.Cells.Item(4, 16).Value = 0.1046219522368387
$ws.Cells.Item(4, 17).Value = 7.349548937773335
$ws.Cells.Item(4, 18).Value = 66.14594043996
$ws.Cells.Item(4, 19).Value = 0.006339773035920536
$ws.Cells.Item(4, 20).Value = 0.006339773035920535

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Clcf1"
$ws.Cells.Item(5, 3).Value = "Lifr"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.6103310000000001
$ws.Cells.Item(5, 8).Value = 1.830993
$ws.Cells.Item(5, 9).Value = 0.06059696746595619
$ws.Cells.Item(5, 10).Value = 0.06059696746595618
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 13.81063466666667
$ws.Cells.Item(5, 14).Value = 41.431904
$ws.Cells.Item(5, 15).Value = 0.1199889353449366
$ws.Cells.Item(5, 16).Value = 0.1199889353449366
$ws.Cells.Item(5, 17).Value = 8.429058466741333
$ws.Cells.Item(5, 18).Value = 75.861526200672
$ws.Cells.Item(5, 19).Value = 0.007270965611371845
$ws.Cells.Item(5, 20).Value = 0.007270965611371845

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Clcf1"
$ws.Cells.Item(6, 3).Value = "Lifr"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 2.690388333333333
$ws.Cells.Item(6, 8).Value = 8.071165
$ws.Cells.Item(6, 9).Value = 0.2671163259047764
$ws.Cells.Item(6, 10).Value = 0.2671163259047764
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 31.44605333333334
$ws.Cells.Item(6, 14).Value = 94.33816
$ws.Cells.Item(6, 15).Value = 0.273208187120734
$ws.Cells.Item(6, 16).Value = 0.273208187120734
$ws.Cells.Item(6, 17).Value = 84.60209501737778
$ws.Cells.Item(6, 18).Value = 761.4188551564001
$ws.Cells.Item(6, 19).Value = 0.07297836715079513
$ws.Cells.Item(6, 20).Value = 0.07297836715079513

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Clcf1"
$ws.Cells.Item(7, 3).Value = "Lifr"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 2.690388333333333
$ws.Cells.Item(7, 8).Value = 8.071165
$ws.Cells.Item(7, 9).Value = 0.2671163259047764
$ws.Cells.Item(7, 10).Value = 0.2671163259047764
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 57.80064033333333
$ws.Cells.Item(7, 14).Value = 173.401921
$ws.Cells.Item(7, 15).Value = 0.5021809252974907
$ws.Cells.Item(7, 16).Value = 0.5021809252974908
$ws.Cells.Item(7, 17).Value = 155.5061684119961
$ws.Cells.Item(7, 18).Value = 1399.555515707965
$ws.Cells.Item(7, 19).Value = 0.1341407237049267
$ws.Cells.Item(7, 20).Value = 0.1341407237049267

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Clcf1"
$ws.Cells.Item(8, 3).Value = "Lifr"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.690388333333333
$ws.Cells.Item(8, 8).Value = 8.071165
$ws.Cells.Item(8, 9).Value = 0.2671163259047764
$ws.Cells.Item(8, 10).Value = 0.2671163259047764
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 12.04190666666667
$ws.Cells.Item(8, 14).Value = 36.12572
$ws.Cells.Item(8, 15).Value = 0.1046219522368387
$ws.Cells.Item(8, 16).Value = 0.1046219522368387
$ws.Cells.Item(8, 17).Value = 32.39740520708889
$ws.Cells.Item(8, 18).Value = 291.5766468638
$ws.Cells.Item(8, 19).Value = 0.02794623149048935
$ws.Cells.Item(8, 20).Value = 0.02794623149048935

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Clcf1"
$ws.Cells.Item(9, 3).Value = "Lifr"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.690388333333333
$ws.Cells.Item(9, 8).Value = 8.071165
$ws.Cells.Item(9, 9).Value = 0.2671163259047764
$ws.Cells.Item(9, 10).Value = 0.2671163259047764
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 13.81063466666667
$ws.Cells.Item(9, 14).Value = 41.431904
$ws.Cells.Item(9, 15).Value = 0.1199889353449366
$ws.Cells.Item(9, 16).Value = 0.1199889353449366
$ws.Cells.Item(9, 17).Value = 37.15597038312889
$ws.Cells.Item(9, 18).Value = 334.40373344816
$ws.Cells.Item(9, 19).Value = 0.03205100355856524
$ws.Cells.Item(9, 20).Value = 0.03205100355856524

# Row 10
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Clcf1"
$ws.Cells.Item(10, 3).Value = "Lifr"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 1.24836
$ws.Cells.Item(10, 8).Value = 3.74508
$ws.Cells.Item(10, 9).Value = 0.1239439423948662
$ws.Cells.Item(10, 10).Value = 0.1239439423948662
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 31.44605333333334
$ws.Cells.Item(10, 14).Value = 94.33816
$ws.Cells.Item(10, 15).Value = 0.273208187120734
$ws.Cells.Item(10, 16).Value = 0.273208187120734
$ws.Cells.Item(10, 17).Value = 39.2559951392
$ws.Cells.Item(10, 18).Value = 353.3039562528
$ws.Cells.Item(10, 19).Value = 0.03386249980629807
$ws.Cells.Item(10, 20).Value = 0.03386249980629807

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Clcf1"
$ws.Cells.Item(11, 3).Value = "Lifr"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 1.24836
$ws.Cells.Item(11, 8).Value = 3.74508
$ws.Cells.Item(11, 9).Value = 0.1239439423948662
$ws.Cells.Item(11, 10).Value = 0.1239439423948662
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 57.80064033333333
$ws.Cells.Item(11, 14).Value = 173.401921
$ws.Cells.Item(11, 15).Value = 0.5021809252974907
$ws.Cells.Item(11, 16).Value = 0.5021809252974908
$ws.Cells.Item(11, 17).Value = 72.15600736651999
$ws.Cells.Item(11, 18).Value = 649.4040662986799
$ws.Cells.Item(11, 19).Value = 0.06224228367687277
$ws.Cells.Item(11, 20).Value = 0.06224228367687278

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Clcf1"
$ws.Cells.Item(12, 3).Value = "Lifr"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 1.24836
$ws.Cells.Item(12, 8).Value = 3.74508
$ws.Cells.Item(12, 9).Value = 0.1239439423948662
$ws.Cells.Item(12, 10).Value = 0.1239439423948662
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 12.04190666666667
$ws.Cells.Item(12, 14).Value = 36.12572
$ws.Cells.Item(12, 15).Value = 0.1046219522368387
$ws.Cells.Item(12, 16).Value = 0.1046219522368387
$ws.Cells.Item(12, 17).Value = 15.0326346064
$ws.Cells.Item(12, 18).Value = 135.2937114576
$ws.Cells.Item(12, 19).Value = 0.01296725722128117
$ws.Cells.Item(12, 20).Value = 0.01296725722128117

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Clcf1"
$ws.Cells.Item(13, 3).Value = "Lifr"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 1.24836
$ws.Cells.Item(13, 8).Value = 3.74508
$ws.Cells.Item(13, 9).Value = 0.1239439423948662
$ws.Cells.Item(13, 10).Value = 0.1239439423948662
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 13.81063466666667
$ws.Cells.Item(13, 14).Value = 41.431904
$ws.Cells.Item(13, 15).Value = 0.1199889353449366
$ws.Cells.Item(13, 16).Value = 0.1199889353449366
$ws.Cells.Item(13, 17).Value = 17.24064389248
$ws.Cells.Item(13, 18).Value = 155.16579503232
$ws.Cells.Item(13, 19).Value = 0.01487190169041414
$ws.Cells.Item(13, 20).Value = 0.01487190169041414

# Row 14
$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Clcf1"
$ws.Cells.Item(14, 3).Value = "Lifr"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 5.522893333333333
$ws.Cells.Item(14, 8).Value = 16.56868
$ws.Cells.Item(14, 9).Value = 0.5483427642344012
$ws.Cells.Item(14, 10).Value = 0.5483427642344012
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 31.44605333333334
$ws.Cells.Item(14, 14).Value = 94.33816
$ws.Cells.Item(14, 15).Value = 0.273208187120734
$ws.Cells.Item(14, 16).Value = 0.273208187120734
$ws.Cells.Item(14, 17).Value = 173.6731983143111
$ws.Cells.Item(14, 18).Value = 1563.0587848288
$ws.Cells.Item(14, 19).Value = 0.1498117325372528
$ws.Cells.Item(14, 20).Value = 0.1498117325372528

# Row 15
$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Clcf1"
$ws.Cells.Item(15, 3).Value = "Lifr"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 5.522893333333333
$ws.Cells.Item(15, 8).Value = 16.56868
$ws.Cells.Item(15, 9).Value = 0.5483427642344012
$ws.Cells.Item(15, 10).Value = 0.5483427642344012
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 57.80064033333333
$ws.Cells.Item(15, 14).Value = 173.401921
$ws.Cells.Item(15, 15).Value = 0.5021809252974907
$ws.Cells.Item(15, 16).Value = 0.5021809252974908
$ws.Cells.Item(15, 17).Value = 319.2267711593644
$ws.Cells.Item(15, 18).Value = 2873.04094043428
$ws.Cells.Item(15, 19).Value = 0.2753672767234154
$ws.Cells.Item(15, 20).Value = 0.2753672767234154

# Row 16
$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Clcf1"
$ws.Cells.Item(16, 3).Value = "Lifr"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 5.522893333333333
$ws.Cells.Item(16, 8).Value = 16.56868
$ws.Cells.Item(16, 9).Value = 0.5483427642344012
$ws.Cells.Item(16, 10).Value = 0.5483427642344012
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 12.04190666666667
$ws.Cells.Item(16, 14).Value = 36.12572
$ws.Cells.Item(16, 15).Value = 0.1046219522368387
$ws.Cells.Item(16, 16).Value = 0.1046219522368387
$ws.Cells.Item(16, 17).Value = 66.50616604995555
$ws.Cells.Item(16, 18).Value = 598.5554944496
$ws.Cells.Item(16, 19).Value = 0.0573686904891476
$ws.Cells.Item(16, 20).Value = 0.0573686904891476

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Clcf1"
$ws.Cells.Item(17, 3).Value = "Lifr"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 5.522893333333333
$ws.Cells.Item(17, 8).Value = 16.56868
$ws.Cells.Item(17, 9).Value = 0.5483427642344012
$ws.Cells.Item(17, 10).Value = 0.5483427642344012
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 13.81063466666667
$ws.Cells.Item(17, 14).Value = 41.431904
$ws.Cells.Item(17, 15).Value = 0.1199889353449366
$ws.Cells.Item(17, 16).Value = 0.1199889353449366
$ws.Cells.Item(17, 17).Value = 76.27466212963554
$ws.Cells.Item(17, 18).Value = 686.4719591667199
$ws.Cells.Item(17, 19).Value = 0.0657950644845854
$ws.Cells.Item(17, 20).Value = 0.0657950644845854
